# Update odds values in row 2 of the active sheet to reflect the latest
# FlashScore data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2"  = 1.5
    "H2"  = 3.8
    "I2"  = 5.7
    "J2"  = 1.98
    "K2"  = 2.25
    "L2"  = 5.5
    "S2"  = 1.3
    "T2"  = 3.32
    "W2"  = 6.4
    "X2"  = 6.5
    "Y2"  = 6.7
    "Z2"  = 9.25
    "AD2" = 6.7
    "AE2" = 12.5
    "AF2" = 45
    "AH2" = 14.5
    "AI2" = 30
    "AJ2" = 14.5
    "AM2" = 37
    "AO2" = 7
    "AP2" = 15
    "AQ2" = 21
    "AR2" = 45
    "AT2" = 3.05
    "AU2" = 7.2
    "AV2" = 55
    "AW2" = 7.4
    "AX2" = 32
    "AY2" = 32
    "BB2" = 350
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
